$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '20.094.55'
$ws.Range('E2').Value = '  -7.71%  '
$ws.Range('D3').Value = '1.428.99'
$ws.Range('E3').Value = '  -7.32%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.001'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.06%  '
$ws.Range('E5').Value = '  +0.05%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '274.79'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -5.07%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.3745'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -3.81%  '
$ws.Range('E8').Value = '  -3.09%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '40.23'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -7.31%  '
$ws.Range('E10').Value = '  -4.78%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.06593'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -8.33%  '
$ws.Range('E12').Value = '  +0.06%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.396'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -4.10%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '17.21'
$ws.Range('D14').Style = 'Normal'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '6.181'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -6.33%  '
$ws.Range('D16').Value = '1.431.31'
$ws.Range('E16').Value = '  -7.17%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.00001011'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -8.67%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '75.61'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -9.27%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.05807'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -11.73%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '1.001'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.08%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '5.690'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -7.35%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '14.52'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -5.53%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '11.12'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +2.46%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.337'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -1.01%  '
$ws.Range('D25').Value = '20.105.56'
$ws.Range('E25').Value = '  -7.73%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.285'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -4.33%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '138.25'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -4.88%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '16.92'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -8.11%  '
$ws.Range('D29').Value = '1.591.11'
$ws.Range('E29').Value = '  -7.28%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '109.41'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -7.08%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.969'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -18.01%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.9093'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -5.97%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.415'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -8.54%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.07777'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -5.20%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '8.411'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -6.16%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '11.40'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +6.86%  '
$ws.Range('B37').Value = 'InternetComputer(DFINITY)'
$ws.Range('C37').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '4.768'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -6.97%  '
$ws.Range('E38').Value = '  +0.05%  '
$ws.Range('B39').Value = 'Hedera'
$ws.Range('C39').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.05695'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -7.13%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.1923'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -5.89%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.119'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -5.43%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.02030'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -8.20%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.299'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -10.14%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.5339'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -7.27%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '3.544'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -5.15%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '12.21'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -7.02%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.5144'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -6.99%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.778'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -5.18%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '109.67'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -6.93%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.053'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -7.12%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.000'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.01%  '
